$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking strings so Excel does not
# convert them to real numbers (which would drop formatting like
# trailing zeros, e.g. "7.20" -> 7.2).
$textCells = @("D5","D6","D7","D10","D15","D16","D19","D20","D21","D22","D25","D27","D29","D30","D33","D34","D35","D36","D38","D39","D42","D43","D45","D48","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "58.847.73"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "3.152.21"
$ws.Range("E3").Value = "  +1.67%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "536.61"
$ws.Range("E5").Value = "  +1.68%  "
$ws.Range("D6").Value = "143.84"
$ws.Range("E6").Value = "  +0.64%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").Value = "3.152.90"
$ws.Range("E8").Value = "  +1.75%  "
$ws.Range("E9").Value = "  +2.26%  "
$ws.Range("D10").Value = "7.20"
$ws.Range("E10").Value = "  -1.93%  "
$ws.Range("E11").Value = "  +0.53%  "
$ws.Range("E12").Value = "  +3.65%  "
$ws.Range("D13").Value = "3.690.96"
$ws.Range("E13").Value = "  +1.57%  "
$ws.Range("E14").Value = "  +3.30%  "
$ws.Range("D15").Value = "25.90"
$ws.Range("E15").Value = "  -3.71%  "
$ws.Range("D16").Value = "0.0000168"
$ws.Range("E16").Value = "  +0.29%  "
$ws.Range("D17").Value = "58.855.43"
$ws.Range("E17").Value = "  +0.53%  "
$ws.Range("D18").Value = "3.157.22"
$ws.Range("E18").Value = "  +1.77%  "
$ws.Range("D19").Value = "6.16"
$ws.Range("E19").Value = "  +0.62%  "
$ws.Range("D20").Value = "12.96"
$ws.Range("E20").Value = "  +0.49%  "
$ws.Range("D21").Value = "8.03"
$ws.Range("E21").Value = "  -0.93%  "
$ws.Range("D22").Value = "343.71"
$ws.Range("E22").Value = "  +0.93%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("E24").Value = "  +1.95%  "
$ws.Range("D25").Value = "67.98"
$ws.Range("E25").Value = "  +2.82%  "
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("D28").Value = "0.0₃0944"
$ws.Range("E28").Value = "  +2.90%  "
$ws.Range("D29").Value = "7.61"
$ws.Range("E29").Value = "  +4.46%  "
$ws.Range("D30").Value = "6.52"
$ws.Range("E30").Value = "  -1.78%  "
$ws.Range("E32").Value = "  +2.01%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "1.22"
$ws.Range("E33").Value = "  +0.59%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "21.21"
$ws.Range("E34").Value = "  +1.12%  "
$ws.Range("D35").Value = "4.84"
$ws.Range("E35").Value = "  +3.11%  "
$ws.Range("D36").Value = "158.30"
$ws.Range("E36").Value = "  +2.77%  "
$ws.Range("E37").Value = "  +3.64%  "
$ws.Range("D38").Value = "26.31"
$ws.Range("E38").Value = "  -2.50%  "
$ws.Range("D39").Value = "1.28"
$ws.Range("E39").Value = "  -2.65%  "
$ws.Range("E40").Value = "  +12.68%  "
$ws.Range("E41").Value = "  -0.54%  "
$ws.Range("D42").Value = "0.712"
$ws.Range("E42").Value = "  +4.97%  "
$ws.Range("D43").Value = "4.04"
$ws.Range("E43").Value = "  +3.80%  "
$ws.Range("D44").Value = "3.189.40"
$ws.Range("E44").Value = "  +1.54%  "
$ws.Range("D45").Value = "36.96"
$ws.Range("E45").Value = "  +0.43%  "
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "2.321.84"
$ws.Range("E47").Value = "  +1.51%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "0.0267"
$ws.Range("E48").Value = "  +3.64%  "
$ws.Range("E49").Value = "  +5.57%  "
$ws.Range("D50").Value = "20.83"
$ws.Range("E50").Value = "  -0.16%  "
$ws.Range("D51").Value = "6.11"
$ws.Range("E51").Value = "  +2.03%  "
